$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.227
$ws.Range("C4").Value = -12.041
$ws.Range("B7").Value = 5.254
$ws.Range("A8").Value = -22.151
$ws.Range("A10").Value = -21.531
$ws.Range("E10").Value = 16.315
$ws.Range("C11").Value = -12.512
$ws.Range("A12").Value = -21.728
$ws.Range("E12").Value = 17.565
$ws.Range("E13").Value = 16.717
$ws.Range("B14").Value = 5.786000000000001
$ws.Range("C14").Value = -12.758
$ws.Range("E14").Value = 16.616
$ws.Range("B15").Value = 5.072
$ws.Range("A18").Value = -21.442
$ws.Range("B18").Value = 7.621
$ws.Range("C18").Value = -11.477
$ws.Range("C19").Value = -11.871
$ws.Range("B20").Value = 6.489999999999999
$ws.Range("C21").Value = -12.016
$ws.Range("A25").Value = -21.896
$ws.Range("C27").Value = -13.419
$ws.Range("B29").Value = 5.42
$ws.Range("E29").Value = 17.18
$ws.Range("B30").Value = 5.552
$ws.Range("B31").Value = 5.522
$ws.Range("C31").Value = -13.16
$ws.Range("E32").Value = 16.613
$ws.Range("B35").Value = 7.304
$ws.Range("E35").Value = 16.521
$ws.Range("A37").Value = -20.727
$ws.Range("C38").Value = -13.055
$ws.Range("B40").Value = 8.627999999999998
$ws.Range("C42").Value = -12.229
$ws.Range("E43").Value = 16.977
$ws.Range("B44").Value = 5.398
$ws.Range("C44").Value = -13.257
$ws.Range("C47").Value = -12.342
$ws.Range("E48").Value = 17.211
$ws.Range("E49").Value = 16.41
$ws.Range("B50").Value = 4.975
$ws.Range("E50").Value = 16.448
$ws.Range("E51").Value = 16.605
$ws.Range("B54").Value = 5.039
$ws.Range("A55").Value = -21.806
$ws.Range("C56").Value = -13.535
$ws.Range("E56").Value = 16.419
$ws.Range("C58").Value = -13.27
$ws.Range("E61").Value = 16.489
$ws.Range("C65").Value = -12.282
$ws.Range("A68").Value = -21.521
$ws.Range("B68").Value = 5.802
$ws.Range("E69").Value = 17.319
$ws.Range("E71").Value = 17.156
$ws.Range("C73").Value = -12.255
$ws.Range("B76").Value = 5.971
$ws.Range("A77").Value = -20.372
$ws.Range("A78").Value = -20.497
$ws.Range("A79").Value = -20.812
$ws.Range("E79").Value = 17.044
$ws.Range("A80").Value = -20.77
$ws.Range("A81").Value = -21.851
$ws.Range("E81").Value = 16.812
$ws.Range("A82").Value = -22.016
$ws.Range("A84").Value = -21.682
$ws.Range("B87").Value = 4.853000000000001
$ws.Range("B88").Value = 5.121
$ws.Range("C90").Value = -13.331
$ws.Range("B92").Value = 6.888
$ws.Range("C92").Value = -11.026
$ws.Range("E92").Value = 17.614
$ws.Range("C94").Value = -10.379
$ws.Range("C95").Value = -11.451
$ws.Range("B96").Value = 6.239999999999999
$ws.Range("B98").Value = 5.44
$ws.Range("A101").Value = -21.073
$ws.Range("B101").Value = 6.294
$ws.Range("C101").Value = -12.194
$ws.Range("A102").Value = -21.32
$ws.Range("B102").Value = 6.241000000000001
